# Add a 2022 column (S) to the sheet, mirroring the R column formatting,
# and update the selection to R17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new year column.
$ws.Range("S4").Value = 2022
$ws.Range("S4").Style = $ws.Range("R4").Style

# Data cells for rows 5-14, all zero, copying the number format/style from
# the corresponding R-column cell.
for ($row = 5; $row -le 14; $row++) {
    $srcCell = $ws.Cells.Item($row, 18)   # column R
    $dstCell = $ws.Cells.Item($row, 19)   # column S
    $dstCell.Value = 0
    $dstCell.Style = $srcCell.Style
}

# Update the active selection.
$ws.Range("R17").Select()
